# Removed InkPath QR code for the live session as advised by the SEDarc team.
# The last slide (sldId 407) only contained the InkPath QR code picture that
# was added for the live session; it is no longer needed, so delete the slide.

$p = $ppt.ActivePresentation

$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq 407) {
        $targetSlide = $candidate
        break
    }
}

if ($targetSlide -eq $null) {
    # Fall back to the last slide if the SlideID lookup ever fails.
    $targetSlide = $p.Slides.Item($p.Slides.Count)
}

$targetSlide.Delete()
